# Apply updated odds values to Sheet1, matching the source XLSX diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.65
$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.3
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.75

# Row 3
$ws.Range("AP3").Value = 29
$ws.Range("AR3").Value = 81
$ws.Range("AS3").Value = 301
$ws.Range("AT3").Value = 2.2
$ws.Range("AU3").Value = 10
$ws.Range("AZ3").Value = 126
$ws.Range("BB3").Value = 501
$ws.Range("G3").Value = 1.85
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 6
$ws.Range("P3").Value = 2.32
$ws.Range("Q3").Value = 2.7
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53

# Row 4
$ws.Range("AQ4").Value = 34
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 6
$ws.Range("S4").Value = 1.57

# Row 5
$ws.Range("S5").Value = 1.4

# Row 7
$ws.Range("AA7").Value = 29
$ws.Range("AB7").Value = 40
$ws.Range("AC7").Value = 6.9
$ws.Range("AF7").Value = 75
$ws.Range("AG7").Value = 700
$ws.Range("AJ7").Value = 9.5
$ws.Range("AK7").Value = 29
$ws.Range("AL7").Value = 23
$ws.Range("AN7").Value = 4.9
$ws.Range("AO7").Value = 16.5
$ws.Range("AP7").Value = 22
$ws.Range("AQ7").Value = 80
$ws.Range("AR7").Value = 110
$ws.Range("AT7").Value = 2.42
$ws.Range("AW7").Value = 4.35
$ws.Range("AX7").Value = 13.5
$ws.Range("G7").Value = 3.05
$ws.Range("H7").Value = 2.75
$ws.Range("I7").Value = 2.52
$ws.Range("J7").Value = 3.55
$ws.Range("K7").Value = 1.95
$ws.Range("L7").Value = 3.05
$ws.Range("N7").Value = 6.75
$ws.Range("P7").Value = 2.52
$ws.Range("T7").Value = 2.45
$ws.Range("U7").Value = 1.8
$ws.Range("V7").Value = 1.8
$ws.Range("X7").Value = 15.5
$ws.Range("Y7").Value = 10.75

# Row 8
$ws.Range("AW8").Value = 6
$ws.Range("AY8").Value = 29
$ws.Range("BB8").Value = 201
$ws.Range("BC8").Value = 126
$ws.Range("BD8").Value = 126
$ws.Range("L8").Value = 4.33
$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.22
$ws.Range("U8").Value = 1.73

# Row 9
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 13
$ws.Range("AD9").Value = 7.5
$ws.Range("AH9").Value = 17
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 17
$ws.Range("AK9").Value = 51
$ws.Range("AO9").Value = 7
$ws.Range("AP9").Value = 15
$ws.Range("AQ9").Value = 21
$ws.Range("AT9").Value = 3.5
$ws.Range("AV9").Value = 41
$ws.Range("AW9").Value = 7.5
$ws.Range("AX9").Value = 29
$ws.Range("BB9").Value = 151
$ws.Range("BD9").Value = 151
$ws.Range("G9").Value = 1.48
$ws.Range("H9").Value = 4.33
$ws.Range("I9").Value = 6.25
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 2.5
$ws.Range("L9").Value = 5.5
$ws.Range("N9").Value = 13
$ws.Range("O9").Value = 1.17
$ws.Range("P9").Value = 4.33
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 2.4
$ws.Range("U9").Value = 1.73
$ws.Range("W9").Value = 8
$ws.Range("X9").Value = 8.5
$ws.Range("Z9").Value = 12

# Row 10
$ws.Range("AW10").Value = 5.5
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.3
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 1.04
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 1.25

# Row 11
$ws.Range("R11").Value = 1.41

# Row 12
$ws.Range("AN12").Value = 4.33
$ws.Range("AO12").Value = 15
$ws.Range("AQ12").Value = 51
$ws.Range("AX12").Value = 19
$ws.Range("G12").Value = 2.35
$ws.Range("I12").Value = 3.1
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 1.5

# Row 13
$ws.Range("AN13").Value = 3.2
$ws.Range("AT13").Value = 2.75
$ws.Range("G13").Value = 1.4
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 9.5
$ws.Range("S13").Value = 1.4
$ws.Range("T13").Value = 2.75
$ws.Range("U13").Value = 2.2
$ws.Range("V13").Value = 1.62
